$d = $word.ActiveDocument

$replacements = @(
    @("2025-05-13 Tuesday", "2025-05-14 Wednesday"),
    @("532÷9=", "988÷6="),
    @("909÷7=", "487÷6="),
    @("471÷7=", "126÷5="),
    @("880÷8=", "917÷3="),
    @("835÷7=", "519÷9="),
    @("327÷2=", "486÷6="),
    @("649÷2=", "683÷9="),
    @("948÷9=", "540÷2="),
    @("214÷2=", "353÷4="),
    @("768÷7=", "539÷4="),
    @("279÷3=", "628÷8="),
    @("712÷8=", "450÷7="),
    @("476÷7=", "464÷5="),
    @("646÷9=", "581÷4="),
    @("530÷3=", "942÷9="),
    @("214÷9=", "796÷8="),
    @("148÷9=", "156÷6="),
    @("750÷6=", "929÷3="),
    @("248÷3=", "235÷3="),
    @("978÷4=", "132÷6="),
    @("914÷5=", "908÷7="),
    @("884÷2=", "243÷6="),
    @("193÷2=", "409÷7="),
    @("276÷3=", "894÷6="),
    @("288÷5=", "149÷7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
